$d = $word.ActiveDocument

# Locate the paragraph that contains the "Ver no Jupiter..." footer text and
# the one that contains the "(c) 2020 ... Creative Commons Attribution" footer
# text, then remove them together with the blank separator paragraph that
# precedes them (the blank line right after the "8800011: Canto Coral III
# (Requisito)" requirement line). The paragraph mark that terminates the
# trailing copyright paragraph is included in the deleted range so the whole
# three-paragraph block collapses away cleanly.

$startPara = $null
$endPara = $null

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    $text = $para.Range.Text

    if ($text -like "*Ver no Jupiter Salvar em pdf Salvar em docx*") {
        # The blank paragraph immediately before this one is the separator
        # that should disappear along with the footer block.
        $startPara = $d.Paragraphs.Item($i - 1)
    }

    if ($text -like "*Powered by Jekyll and Github pages*") {
        $endPara = $para
    }
}

if ($startPara -ne $null -and $endPara -ne $null) {
    $deleteRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $deleteRange.Delete()
}
